# Weekly fruit/vegetable update: insert a new price-report row at row 350
# (pushing the existing rows 350-380 down to 351-381) and populate the new
# row with this week's data for Feria Lagunitas de Puerto Montt - Brócoli.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 350, shifting rows 350:380
# down to 351:381 (and the sheet dimension grows from R380 to R381).
$ws.Rows("350:350").Insert()

# Populate the newly inserted row 350 with the new weekly record.
$ws.Range("A350").Value2 = 4
$ws.Range("B350").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C350").Value2 = "Los Lagos"
$ws.Range("D350").Value2 = 44783
$ws.Range("E350").Value2 = 10
$ws.Range("F350").Value2 = 100112023
$ws.Range("G350").Value2 = "Brócoli"
$ws.Range("H350").Value2 = "Sin especificar"
$ws.Range("I350").Value2 = "Primera"
$ws.Range("J350").Value2 = 100
$ws.Range("K350").Value2 = 1500
$ws.Range("L350").Value2 = 1500
$ws.Range("M350").Value2 = 1500
$ws.Range("N350").Value2 = "`$/unidad"
$ws.Range("O350").Value2 = "Región Metropolitana"
$ws.Range("P350").Value2 = 1500
$ws.Range("Q350").Value2 = 1
$ws.Range("R350").Value2 = "Hortaliza"
